$wb = $excel.ActiveWorkbook
$wsInfo = $wb.Worksheets.Item("info")
$wsData = $wb.Worksheets.Item("data")

# --- Sheet "info": update the Last update timestamp ---
$wsInfo.Range("B2").Value = "2021-09-24 11:16:54"

# --- Sheet "data": update summary label text (column G) ---
$gLabels = @{
    "G2"  = "1190 Gt CO2"
    "G3"  = "1190 Gt CO2"
    "G4"  = "1190 Gt CO2"
    "G5"  = "1190 Gt CO2"
    "G6"  = "1190 Gt CO2"
    "G7"  = "1190 Gt CO2"
    "G38" = "605 Gt CO2"
    "G39" = "605 Gt CO2"
    "G40" = "605 Gt CO2"
    "G41" = "605 Gt CO2"
    "G42" = "605 Gt CO2"
    "G43" = "605 Gt CO2"
    "G50" = "51 Gt CO2eq"
    "G51" = "51 Gt CO2eq"
    "G52" = "51 Gt CO2eq"
    "G53" = "51 Gt CO2eq"
    "G54" = "51 Gt CO2eq"
    "G55" = "51 Gt CO2eq"
    "G56" = "58 Gt CO2eq"
    "G57" = "58 Gt CO2eq"
    "G58" = "58 Gt CO2eq"
    "G59" = "58 Gt CO2eq"
    "G60" = "58 Gt CO2eq"
    "G61" = "58 Gt CO2eq"
}
foreach ($ref in $gLabels.Keys) {
    $wsData.Range($ref).Value = $gLabels[$ref]
}

# --- Sheet "data": update numeric columns B, E, F with the final data version ---
$numericUpdates = @{
    "B2" = 41241446843.9533
    "E2" = 1189679505864.58
    "F2" = 0.0346660143682829
    "B3" = 326712801404.424
    "E3" = 1189679505864.58
    "F3" = 0.274622534719542
    "B4" = 599201738199.917
    "E4" = 1189679505864.58
    "F4" = 0.503666521316138
    "B5" = 114687068188.585
    "E5" = 1189679505864.58
    "F5" = 0.0964016507162051
    "B6" = 57729320270.8385
    "E6" = 1189679505864.58
    "F6" = 0.0485251027577251
    "B7" = 50107130956.8652
    "E7" = 1189679505864.58
    "F7" = 0.0421181761221066
    "B8" = 1494770229.59111
    "E8" = 36521769730.4009
    "F8" = 0.0409281981849542
    "B9" = 17356975844.8493
    "E9" = 36521769730.4009
    "F9" = 0.475250130893883
    "B10" = 11245656786.8144
    "E10" = 36521769730.4009
    "F10" = 0.307916534982517
    "B11" = 2463828549.26231
    "E11" = 36521769730.4009
    "F11" = 0.0674619156587971
    "B12" = 1770512169.96196
    "E12" = 36521769730.4009
    "F12" = 0.0484782687978064
    "B13" = 2190026149.92183
    "E13" = 36521769730.4009
    "F13" = 0.0599649514820428
    "B26" = 1302819380
    "E26" = 7618629943
    "F26" = 0.171004418083993
    "B27" = 3970583353
    "E27" = 7618629943
    "F27" = 0.521167635481255
    "B28" = 1192863286
    "E28" = 7618629943
    "F28" = 0.15657188955555
    "B29" = 249086331
    "E29" = 7618629943
    "F29" = 0.0326943732486785
    "B30" = 646430786
    "E30" = 7618629943
    "F30" = 0.0848486920662082
    "B31" = 256846807
    "E31" = 7618629943
    "F31" = 0.0337129915643154
    "B32" = 1842565311.82348
    "E32" = 490738295127.81
    "F32" = 0.00375468010162849
    "B33" = 94976534825.3931
    "E33" = 490738295127.81
    "F33" = 0.193538054332314
    "B34" = 370041679034.565
    "E34" = 490738295127.81
    "F34" = 0.754050952836664
    "B35" = 57825944.5643504
    "E35" = 490738295127.81
    "F35" = 0.000117834587474552
    "B36" = 11944974298.2492
    "E36" = 490738295127.81
    "F36" = 0.0243408236464167
    "B37" = 11874715713.2148
    "E37" = 490738295127.81
    "F37" = 0.0241976544955027
    "B38" = 10293558079.3937
    "E38" = 604680080872.7
    "F38" = 0.0170231472889558
    "B39" = 192043417221.044
    "E39" = 604680080872.7
    "F39" = 0.317595077621672
    "B40" = 227105692972.81
    "E40" = 604680080872.7
    "F40" = 0.375579914332619
    "B41" = 102898575536.138
    "E41" = 604680080872.7
    "F41" = 0.170170274813138
    "B42" = 40225847241.0209
    "E42" = 604680080872.7
    "F42" = 0.0665241811553726
    "B43" = 32112989822.2935
    "E43" = 604680080872.7
    "F43" = 0.0531074047882421
    "B44" = 29105323452.7362
    "E44" = 94261129864.0724
    "F44" = 0.308773335251837
    "B45" = 39692849357.9866
    "E45" = 94261129864.0724
    "F45" = 0.421094563742499
    "B46" = 2054366192.54176
    "E46" = 94261129864.0724
    "F46" = 0.0217944151051894
    "B47" = 11730666707.8825
    "E47" = 94261129864.0724
    "F47" = 0.124448611265307
    "B48" = 5558498731.56837
    "E48" = 94261129864.0724
    "F48" = 0.0589691502699353
    "B49" = 6119425421.35697
    "E49" = 94261129864.0724
    "F49" = 0.0649199243652329
    "B50" = 3433150433.67798
    "E50" = 51126029486.308
    "F50" = 0.0671507345313685
    "B51" = 23078805804.8676
    "E51" = 51126029486.308
    "F51" = 0.45141009455953
    "B52" = 14284243064.0032
    "E52" = 51126029486.308
    "F52" = 0.279392771305049
    "B53" = 3416156396.33553
    "E53" = 51126029486.308
    "F53" = 0.0668183395162813
    "B54" = 3695402314.40367
    "E54" = 51126029486.308
    "F54" = 0.072280252378944
    "B55" = 3218271473.01995
    "E55" = 51126029486.308
    "F55" = 0.0629478077088273
    "B56" = 5209052921.67265
    "E56" = 57731068333.908
    "F56" = 0.0902296297644149
    "B57" = 25718302637.0276
    "E57" = 57731068333.908
    "F57" = 0.445484613038455
    "B58" = 14143018287.6832
    "E58" = 57731068333.908
    "F58" = 0.244981059520363
    "B59" = 3505445742.36753
    "E59" = 57731068333.908
    "F59" = 0.0607202645565565
    "B60" = 5938507787.225
    "E60" = 57731068333.908
    "F60" = 0.102865024996204
    "B61" = 3216740957.93195
    "E61" = 57731068333.908
    "F61" = 0.0557194081240069
}
foreach ($ref in $numericUpdates.Keys) {
    $wsData.Range($ref).Value = $numericUpdates[$ref]
}

Write-Output "done"